$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.015.35"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").Value = "1.638.68"
$ws.Range("E3").Value = "  -1.60%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.015"
$cell.ClearFormats()
$ws.Range("E4").Value = "  +0.78%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "215.90"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -1.14%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.016"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.79%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5000"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -2.77%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2577"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +0.20%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06423"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -0.31%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "19.48"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -2.31%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07756"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +1.28%  "

$ws.Range("D12").Value = "1.644.18"
$ws.Range("E12").Value = "  -1.41%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.251"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -1.89%  "

$ws.Range("D14").Value = "1.865.54"
$ws.Range("E14").Value = "  -1.53%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.5449"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -1.55%  "

$ws.Range("D16").Value = "0.0₅7940"
$ws.Range("E16").Value = "  -1.14%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "63.67"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("D18").Value = "26.034.09"
$ws.Range("E18").Value = "  -1.49%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "1.016"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +0.79%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "204.85"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -2.23%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.311"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -2.37%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.01"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -0.85%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.971"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +1.67%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.016"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +0.83%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.972"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +13.57%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "141.80"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -2.27%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1154"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -0.76%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "15.78"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +0.16%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "6.802"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -2.73%  "

$ws.Range("E30").Value = "  -3.57%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.244"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -1.39%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.266"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -2.99%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.203"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -0.30%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.544"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -1.95%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.357"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -0.77%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.8924"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -3.68%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.613"
$cell.ClearFormats()

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.5648"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("D39").Value = "1.124.80"
$ws.Range("E39").Value = "  -2.40%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.01565"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -1.99%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "2.585"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +0.72%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.016"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +0.85%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.639"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -0.02%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.8175"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -2.99%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "99.84"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("D46").Value = "1.774.85"
$ws.Range("E46").Value = "  -1.72%  "

$ws.Range("D47").Value = "0.0₈113"
$ws.Range("E47").Value = "  +1.02%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.4550"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +1.23%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.019"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +1.41%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "54.83"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -1.94%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.05039"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -1.45%  "
